# Updated cryptos list on Sun Mar 10 06:19:07 UTC 2024 with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-Text($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
}

# Unicode subscript digits used in very-small-price notation (e.g. 0.0\u2083...)
$sub3 = [char]0x2083
$sub6 = [char]0x2086

# Row 2 - Bitcoin
Set-Text "D2" "69.344.75"
Set-Text "E2" "  +1.68%  "

# Row 3 - Ethereum
Set-Text "D3" "3.937.96"
Set-Text "E3" "  +0.42%  "

# Row 4 - TetherUSD
Set-Text "E4" "  +0.11%  "

# Row 5 - BNB
Set-Text "D5" "493.46"

# Row 6 - Solana
Set-Text "D6" "147.55"
Set-Text "E6" "  -0.60%  "

# Row 7 - XRP
Set-Text "D7" "0.623"
Set-Text "E7" "  -1.00%  "

# Row 8 - USDC
Set-Text "E8" "  +0.01%  "

# Row 9 - Cardano
Set-Text "D9" "0.734"
Set-Text "E9" "  +0.26%  "

# Row 10 - Dogecoin
Set-Text "E10" "  +4.23%  "

# Row 11 - ShibaInu
Set-Text "E11" "  -0.98%  "

# Row 12
Set-Text "D12" "43.37"
Set-Text "E12" "  +0.80%  "

# Row 13
Set-Text "D13" "10.44"
Set-Text "E13" "  -1.85%  "

# Row 14
Set-Text "D14" "4.568.50"
Set-Text "E14" "  +0.49%  "

# Row 15
Set-Text "D15" "3.945.07"
Set-Text "E15" "  +0.31%  "

# Row 16
Set-Text "D16" "14.25"
Set-Text "E16" "  -3.40%  "

# Row 17 - TRON
Set-Text "E17" "  -0.75%  "

# Row 18 - Polygon
Set-Text "E18" "  +3.94%  "

# Row 19
Set-Text "D19" "19.89"
Set-Text "E19" "  -0.47%  "

# Row 20
Set-Text "D20" "69.300.53"
Set-Text "E20" "  +1.44%  "

# Row 21
Set-Text "D21" "440.10"
Set-Text "E21" "  -0.43%  "

# Row 22
Set-Text "D22" "3.46"
Set-Text "E22" "  +1.23%  "

# Row 23
Set-Text "D23" "14.52"
Set-Text "E23" "  -2.52%  "

# Row 24
Set-Text "D24" "88.70"
Set-Text "E24" "  +0.10%  "

# Row 25
Set-Text "D25" "12.11"
Set-Text "E25" "  +7.44%  "

# Row 26
Set-Text "D26" "3.79"
Set-Text "E26" "  +4.34%  "

# Row 27 - Filecoin
Set-Text "E27" "  -2.95%  "

# Row 28
Set-Text "D28" "37.17"
Set-Text "E28" "  -4.18%  "

# Row 29
Set-Text "D29" "5.65"
Set-Text "E29" "  -3.78%  "

# Row 30 - Bittensor
Set-Text "D30" "705.58"
Set-Text "E30" "  -2.03%  "

# Row 31 - becomes Hedera (was Cosmos)
Set-Text "B31" "Hedera"
Set-Text "C31" "https://coinranking.com/coin/jad286TjB+hedera-hbar"
Set-Text "D31" "0.131"
Set-Text "E31" "  +0.11%  "

# Row 32 - becomes Cosmos (was Hedera)
Set-Text "B32" "Cosmos"
Set-Text "C32" "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
Set-Text "D32" "13.36"
Set-Text "E32" "  -1.24%  "

# Row 33 - Toncoin
Set-Text "E33" "  +0.36%  "

# Row 34 - TheGraph
Set-Text "D34" "0.473"
Set-Text "E34" "  +16.50%  "

# Row 35 - PEPE
Set-Text "D35" "0.0${sub3}0908"
Set-Text "E35" "  -0.76%  "

# Row 36 - becomes OKB (was NEARProtocol)
Set-Text "B36" "OKB"
Set-Text "C36" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-Text "D36" "61.68"
Set-Text "E36" "  +3.79%  "

# Row 37 - becomes NEARProtocol (was OKB)
Set-Text "B37" "NEARProtocol"
Set-Text "C37" "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-Text "D37" "6.07"
Set-Text "E37" "  +0.78%  "

# Row 38 - InjectiveProtocol
Set-Text "D38" "40.71"
Set-Text "E38" "  -2.49%  "

# Row 39
Set-Text "D39" "0.151"
Set-Text "E39" "  +0.61%  "

# Row 40
Set-Text "D40" "0.998"
Set-Text "E40" "  -0.34%  "

# Row 41 - FirstDigitalUSD
Set-Text "E41" "  +0.03%  "

# Row 42
Set-Text "D42" "0.0490"
Set-Text "E42" "  +1.60%  "

# Row 43
Set-Text "D43" "2.91"
Set-Text "E43" "  -1.47%  "

# Row 44 - ThetaToken
Set-Text "E44" "  -3.43%  "

# Row 45
Set-Text "D45" "3.00"
Set-Text "E45" "  +2.31%  "

# Row 46 - Stellar
Set-Text "E46" "  +0.37%  "

# Row 47
Set-Text "D47" "3.38"
Set-Text "E47" "  +7.86%  "

# Row 48
Set-Text "D48" "0.0${sub6}0358"
Set-Text "E48" "  +0.69%  "

# Row 49
Set-Text "D49" "2.99"
Set-Text "E49" "  +5.29%  "

# Row 50 - LidoDAOToken
Set-Text "E50" "  -1.13%  "

# Row 51
Set-Text "D51" "144.25"
Set-Text "E51" "  -0.78%  "

Write-Output "applied"
